$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The rule name for the last row (row 11, column B) changes from "R40" to
# the text "1". Assigning a plain numeric-looking string via .Value would
# store it as a number (and changing NumberFormat to force text would
# allocate a brand new cell style). Instead, produce the text through a
# formula and then convert it to a literal value in place, which keeps the
# cell's existing style/format untouched.
$cell = $ws.Range("B11")
$cell.Formula = '=TEXT(1,"0")'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
